# Append 7 new species-observation records (rows 34-40) to the "Artfynd"
# sheet, mirroring the shape of the existing rows (same populated columns,
# same blank/text conventions).

$rows = @(
  [PSCustomObject]@{ Row = 34; A = 112038596; B = 90087; C = 'Ovaliderad'; D = 'LC'; E = 3298; F = 'Trådticka'; G = 'Climacocystis borealis'; H = '(Fr.) Kotl. & Pouzar'; I = $null; P = 'Västanå, Mpd'; Q = 616076.0611235843; R = 6895427.595461337; S = 5; T = 'Västernorrland'; U = 'Sundsvall'; V = 'Medelpad'; W = 'Njurunda'; Y = '2023-09-09'; Z = '00:00'; AA = '2023-09-09'; AB = '00:00'; AD = $false; AE = $false; AG = $false; AT = $null; AW = 'Olle Finnström'; AX = 'Olle Finnström'; AY = $null }
  [PSCustomObject]@{ Row = 35; A = 112038601; B = 73634; C = 'Ovaliderad'; D = 'LC'; E = 6426; F = 'Kattfotslav'; G = 'Felipes leucopellaeus'; H = '(Ach.) Frisch & G.Thor'; I = $null; P = 'Västanå, Mpd'; Q = 616012.5978259755; R = 6895611.944218947; S = 5; T = 'Västernorrland'; U = 'Sundsvall'; V = 'Medelpad'; W = 'Njurunda'; Y = '2023-09-09'; Z = '00:00'; AA = '2023-09-09'; AB = '00:00'; AD = $false; AE = $false; AG = $false; AT = $null; AW = 'Olle Finnström'; AX = 'Olle Finnström'; AY = $null }
  [PSCustomObject]@{ Row = 36; A = 112038600; B = 86223; C = 'Ovaliderad'; D = 'NT'; E = 4412; F = 'Äggvaxskivling'; G = 'Hygrophorus karstenii'; H = 'Sacc. & Cub.'; I = $null; P = 'Västanå, Mpd'; Q = 616034.1211971109; R = 6895585.10294092; S = 5; T = 'Västernorrland'; U = 'Sundsvall'; V = 'Medelpad'; W = 'Njurunda'; Y = '2023-09-09'; Z = '00:00'; AA = '2023-09-09'; AB = '00:00'; AD = $false; AE = $false; AG = $false; AT = $null; AW = 'Olle Finnström'; AX = 'Olle Finnström'; AY = $null }
  [PSCustomObject]@{ Row = 37; A = 112038604; B = 89845; C = 'Ovaliderad'; D = 'VU'; E = 1209; F = 'Rynkskinn'; G = 'Phlebia centrifuga'; H = 'P.Karst.'; I = $null; P = 'Västanå, Mpd'; Q = 615977.7276359925; R = 6895550.438170813; S = 5; T = 'Västernorrland'; U = 'Sundsvall'; V = 'Medelpad'; W = 'Njurunda'; Y = '2023-09-09'; Z = '00:00'; AA = '2023-09-09'; AB = '00:00'; AD = $false; AE = $false; AG = $false; AT = $null; AW = 'Olle Finnström'; AX = 'Olle Finnström'; AY = $null }
  [PSCustomObject]@{ Row = 38; A = 112038599; B = 89423; C = 'Ovaliderad'; D = 'NT'; E = 5432; F = 'Granticka'; G = 'Porodaedalea chrysoloma'; H = '(Fr.) Fiasson & Niemelä'; I = $null; P = 'Västanå, Mpd'; Q = 616070.2961488151; R = 6895499.860901954; S = 5; T = 'Västernorrland'; U = 'Sundsvall'; V = 'Medelpad'; W = 'Njurunda'; Y = '2023-09-09'; Z = '00:00'; AA = '2023-09-09'; AB = '00:00'; AD = $false; AE = $false; AG = $false; AT = $null; AW = 'Olle Finnström'; AX = 'Olle Finnström'; AY = $null }
  [PSCustomObject]@{ Row = 39; A = 112038602; B = 86223; C = 'Ovaliderad'; D = 'NT'; E = 4412; F = 'Äggvaxskivling'; G = 'Hygrophorus karstenii'; H = 'Sacc. & Cub.'; I = $null; P = 'Västanå, Mpd'; Q = 616026.2967975155; R = 6895553.979090866; S = 5; T = 'Västernorrland'; U = 'Sundsvall'; V = 'Medelpad'; W = 'Njurunda'; Y = '2023-09-09'; Z = '00:00'; AA = '2023-09-09'; AB = '00:00'; AD = $false; AE = $false; AG = $false; AT = $null; AW = 'Olle Finnström'; AX = 'Olle Finnström'; AY = $null }
  [PSCustomObject]@{ Row = 40; A = 112038603; B = 89369; C = 'Ovaliderad'; D = 'LC'; E = 5447; F = 'Vedticka'; G = 'Fuscoporia viticola'; H = '(Schwein.) Murrill'; I = $null; P = 'Västanå, Mpd'; Q = 615968.1934313668; R = 6895405.650930508; S = 5; T = 'Västernorrland'; U = 'Sundsvall'; V = 'Medelpad'; W = 'Njurunda'; Y = '2023-09-09'; Z = '00:00'; AA = '2023-09-09'; AB = '00:00'; AD = $false; AE = $false; AG = $false; AT = $null; AW = 'Olle Finnström'; AX = 'Olle Finnström'; AY = $null }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericCols = @("A","B","E","Q","R","S")
$boolCols = @("AD","AE","AG")
$textAsEnteredCols = @("Y","Z","AA","AB")  # date/time-looking values that must stay as literal text
$plainTextCols = @("C","D","F","G","H","P","T","U","V","W","AW","AX")
$blankTextCols = @("I","AT","AY")          # originally empty text cells

foreach ($row in $rows) {
    $r = $row.Row

    foreach ($col in $numericCols) {
        $ws.Range("$col$r").Value = $row.$col
    }

    foreach ($col in $boolCols) {
        $ws.Range("$col$r").Value = $row.$col
    }

    foreach ($col in $plainTextCols) {
        $ws.Range("$col$r").Value = $row.$col
    }

    foreach ($col in $textAsEnteredCols) {
        # Leading apostrophe forces Excel to store the value as literal text
        # instead of auto-converting the date/time-shaped string to a serial number.
        $ws.Range("$col$r").Value = "'" + $row.$col
        $ws.Range("$col$r").Style = "Normal"
    }

    foreach ($col in $blankTextCols) {
        # Leading apostrophe with nothing after it yields a real empty-string
        # text cell (matches the source file's empty inlineStr cells) rather
        # than clearing the cell back to blank/unset.
        $ws.Range("$col$r").Value = "'"
        $ws.Range("$col$r").Style = "Normal"
    }
}
